# Prepend "stim/" to the image filenames stored in column C (rows 2-9),
# matching the fix described in the commit message: the resource folder
# paths were missing the "stim/" prefix.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 9; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = "stim/" + $cell.Value2
}

# Update the active selection to C9, matching the saved sheet view state.
$ws.Range("C9").Select()
